$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two obsolete data rows (old rows 12-13) so the table shrinks to 10 data rows
$ws.Rows("12:13").Delete()

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Fgf13"
$ws.Cells.Item(2,3).Value = "Scn5a"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.5
$ws.Cells.Item(2,7).Value = 0.1113035
$ws.Cells.Item(2,8).Value = 0.222607
$ws.Cells.Item(2,9).Value = 0.0289261410371415
$ws.Cells.Item(2,10).Value = 0.0273813389619527
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 0.5100455
$ws.Cells.Item(2,14).Value = 1.020091
$ws.Cells.Item(2,15).Value = 0.07988759408568151
$ws.Cells.Item(2,16).Value = 0.07988759408568151
$ws.Cells.Item(2,17).Value = 0.05676984930925
$ws.Cells.Item(2,18).Value = 0.227079397237
$ws.Cells.Item(2,19).Value = 0.002310839813640335
$ws.Cells.Item(2,20).Value = 0.002187429292514933

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Fgf13"
$ws.Cells.Item(3,3).Value = "Scn5a"
$ws.Cells.Item(3,4).Value = "MuSCs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.5
$ws.Cells.Item(3,7).Value = 0.1113035
$ws.Cells.Item(3,8).Value = 0.222607
$ws.Cells.Item(3,9).Value = 0.0289261410371415
$ws.Cells.Item(3,10).Value = 0.0273813389619527
$ws.Cells.Item(3,11).Value = 2
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 5.874494
$ws.Cells.Item(3,14).Value = 11.748988
$ws.Cells.Item(3,15).Value = 0.9201124059143184
$ws.Cells.Item(3,16).Value = 0.9201124059143184
$ws.Cells.Item(3,17).Value = 0.6538517429290001
$ws.Cells.Item(3,18).Value = 2.615406971716
$ws.Cells.Item(3,19).Value = 0.02661530122350117
$ws.Cells.Item(3,20).Value = 0.02519390966943776

# Row 4
$ws.Cells.Item(4,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(4,2).Value = "Fgf13"
$ws.Cells.Item(4,3).Value = "Scn5a"
$ws.Cells.Item(4,4).Value = "ECs"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.406345
$ws.Cells.Item(4,8).Value = 1.219035
$ws.Cells.Item(4,9).Value = 0.1056030832789379
$ws.Cells.Item(4,10).Value = 0.1499450176386367
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 0.5100455
$ws.Cells.Item(4,14).Value = 1.020091
$ws.Cells.Item(4,15).Value = 0.07988759408568151
$ws.Cells.Item(4,16).Value = 0.07988759408568151
$ws.Cells.Item(4,17).Value = 0.2072544386975
$ws.Cells.Item(4,18).Value = 1.243526632185
$ws.Cells.Item(4,19).Value = 0.008436376251184211
$ws.Cells.Item(4,20).Value = 0.01197874670428577

# Row 5
$ws.Cells.Item(5,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(5,2).Value = "Fgf13"
$ws.Cells.Item(5,3).Value = "Scn5a"
$ws.Cells.Item(5,4).Value = "MuSCs"
$ws.Cells.Item(5,5).Value = 1
$ws.Cells.Item(5,6).Value = 0.3333333333333333
$ws.Cells.Item(5,7).Value = 0.406345
$ws.Cells.Item(5,8).Value = 1.219035
$ws.Cells.Item(5,9).Value = 0.1056030832789379
$ws.Cells.Item(5,10).Value = 0.1499450176386367
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 5.874494
$ws.Cells.Item(5,14).Value = 11.748988
$ws.Cells.Item(5,15).Value = 0.9201124059143184
$ws.Cells.Item(5,16).Value = 0.9201124059143184
$ws.Cells.Item(5,17).Value = 2.38707126443
$ws.Cells.Item(5,18).Value = 14.32242758658
$ws.Cells.Item(5,19).Value = 0.09716670702775368
$ws.Cells.Item(5,20).Value = 0.137966270934351

# Row 6
$ws.Cells.Item(6,1).Value = "MuSCs"
$ws.Cells.Item(6,2).Value = "Fgf13"
$ws.Cells.Item(6,3).Value = "Scn5a"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 2
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 3.3023715
$ws.Cells.Item(6,8).Value = 6.604743
$ws.Cells.Item(6,9).Value = 0.8582377352557337
$ws.Cells.Item(6,10).Value = 0.8124035041107618
$ws.Cells.Item(6,11).Value = 2
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 0.5100455
$ws.Cells.Item(6,14).Value = 1.020091
$ws.Cells.Item(6,15).Value = 0.07988759408568151
$ws.Cells.Item(6,16).Value = 0.07988759408568151
$ws.Cells.Item(6,17).Value = 1.68435972290325
$ws.Cells.Item(6,18).Value = 6.737438891613
$ws.Cells.Item(6,19).Value = 0.06856254782312464
$ws.Cells.Item(6,20).Value = 0.06490096137018583

# Row 7
$ws.Cells.Item(7,1).Value = "MuSCs"
$ws.Cells.Item(7,2).Value = "Fgf13"
$ws.Cells.Item(7,3).Value = "Scn5a"
$ws.Cells.Item(7,4).Value = "MuSCs"
$ws.Cells.Item(7,5).Value = 2
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 3.3023715
$ws.Cells.Item(7,8).Value = 6.604743
$ws.Cells.Item(7,9).Value = 0.8582377352557337
$ws.Cells.Item(7,10).Value = 0.8124035041107618
$ws.Cells.Item(7,11).Value = 2
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 5.874494
$ws.Cells.Item(7,14).Value = 11.748988
$ws.Cells.Item(7,15).Value = 0.9201124059143184
$ws.Cells.Item(7,16).Value = 0.9201124059143184
$ws.Cells.Item(7,17).Value = 19.399761562521
$ws.Cells.Item(7,18).Value = 77.599046250084
$ws.Cells.Item(7,19).Value = 0.789675187432609
$ws.Cells.Item(7,20).Value = 0.7475025427405759

# Row 8
$ws.Cells.Item(8,1).Value = "Neutrophils"
$ws.Cells.Item(8,2).Value = "Fgf13"
$ws.Cells.Item(8,3).Value = "Scn5a"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 1
$ws.Cells.Item(8,6).Value = 0.3333333333333333
$ws.Cells.Item(8,7).Value = 0.0004946666666666667
$ws.Cells.Item(8,8).Value = 0.001484
$ws.Cells.Item(8,9).Value = 0.0001285565841718604
$ws.Cells.Item(8,10).Value = 0.0001825365196042254
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 0.5100455
$ws.Cells.Item(8,14).Value = 1.020091
$ws.Cells.Item(8,15).Value = 0.07988759408568151
$ws.Cells.Item(8,16).Value = 0.07988759408568151
$ws.Cells.Item(8,17).Value = 0.0002523025073333334
$ws.Cells.Item(8,18).Value = 0.001513815044
$ws.Cells.Item(8,19).Value = 0.00001027007621336333
$ws.Cells.Item(8,20).Value = 0.0000145824033839554

# Row 9
$ws.Cells.Item(9,1).Value = "Neutrophils"
$ws.Cells.Item(9,2).Value = "Fgf13"
$ws.Cells.Item(9,3).Value = "Scn5a"
$ws.Cells.Item(9,4).Value = "MuSCs"
$ws.Cells.Item(9,5).Value = 1
$ws.Cells.Item(9,6).Value = 0.3333333333333333
$ws.Cells.Item(9,7).Value = 0.0004946666666666667
$ws.Cells.Item(9,8).Value = 0.001484
$ws.Cells.Item(9,9).Value = 0.0001285565841718604
$ws.Cells.Item(9,10).Value = 0.0001825365196042254
$ws.Cells.Item(9,11).Value = 2
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 5.874494
$ws.Cells.Item(9,14).Value = 11.748988
$ws.Cells.Item(9,15).Value = 0.9201124059143184
$ws.Cells.Item(9,16).Value = 0.9201124059143184
$ws.Cells.Item(9,17).Value = 0.002905916365333334
$ws.Cells.Item(9,18).Value = 0.017435498192
$ws.Cells.Item(9,19).Value = 0.0001182865079584971
$ws.Cells.Item(9,20).Value = 0.00016795411622027

# Row 10
$ws.Cells.Item(10,1).Value = "Resolving-Mac"
$ws.Cells.Item(10,2).Value = "Fgf13"
$ws.Cells.Item(10,3).Value = "Scn5a"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 1
$ws.Cells.Item(10,6).Value = 0.3333333333333333
$ws.Cells.Item(10,7).Value = 0.027337
$ws.Cells.Item(10,8).Value = 0.082011
$ws.Cells.Item(10,9).Value = 0.007104483844015123
$ws.Cells.Item(10,10).Value = 0.01008760276904456
$ws.Cells.Item(10,11).Value = 2
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 0.5100455
$ws.Cells.Item(10,14).Value = 1.020091
$ws.Cells.Item(10,15).Value = 0.07988759408568151
$ws.Cells.Item(10,16).Value = 0.07988759408568151
$ws.Cells.Item(10,17).Value = 0.0139431138335
$ws.Cells.Item(10,18).Value = 0.083658683001
$ws.Cells.Item(10,19).Value = 0.0005675601215189624
$ws.Cells.Item(10,20).Value = 0.0008058743153110287

# Row 11
$ws.Cells.Item(11,1).Value = "Resolving-Mac"
$ws.Cells.Item(11,2).Value = "Fgf13"
$ws.Cells.Item(11,3).Value = "Scn5a"
$ws.Cells.Item(11,4).Value = "MuSCs"
$ws.Cells.Item(11,5).Value = 1
$ws.Cells.Item(11,6).Value = 0.3333333333333333
$ws.Cells.Item(11,7).Value = 0.027337
$ws.Cells.Item(11,8).Value = 0.082011
$ws.Cells.Item(11,9).Value = 0.007104483844015123
$ws.Cells.Item(11,10).Value = 0.01008760276904456
$ws.Cells.Item(11,11).Value = 2
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 5.874494
$ws.Cells.Item(11,14).Value = 11.748988
$ws.Cells.Item(11,15).Value = 0.9201124059143184
$ws.Cells.Item(11,16).Value = 0.9201124059143184
$ws.Cells.Item(11,17).Value = 0.160591042478
$ws.Cells.Item(11,18).Value = 0.9635462548680001
$ws.Cells.Item(11,19).Value = 0.00653692372249616
$ws.Cells.Item(11,20).Value = 0.009281728453733533
